# Finished Field Download component.
# - Rename the sheet to its final name.
# - Move the active-cell selection from D5 to D8.
# - Drop the redundant per-row "custom format" (s="3" customFormat="1")
#   that every data row carried, while preserving each cell's own style
#   (bold header / left alignment / wrap-text) and the explicit row
#   heights used for the wrapped description cells.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ClearFormats() strips the row-level default style (the source of the
# "customFormat" row attribute) without touching the explicit row
# heights already set on rows 5, 7, 8, 9, 10, 11, 12.
for ($r = 1; $r -le 12; $r++) {
    $ws.Rows.Item($r).ClearFormats()
}

# Restore the visual formatting that the original per-cell styles
# carried (vertical-top alignment everywhere, bold header row, left
# alignment for the Field Name column, and wrapped Description cells).
$ws.Range("A1:D12").VerticalAlignment = -4160   # xlTop

$ws.Range("A1:D1").Font.Bold = $true
$ws.Range("B1:B12").HorizontalAlignment = -4131 # xlLeft

$ws.Range("D5").WrapText = $true
$ws.Range("D7").WrapText = $true
$ws.Range("D8").WrapText = $true
$ws.Range("D9").WrapText = $true
$ws.Range("D10").WrapText = $true
$ws.Range("D11").WrapText = $true
$ws.Range("D12").WrapText = $true

# Rename the worksheet tab.
$ws.Name = "food_event_fields"

# Move the selection to D8.
$ws.Range("D8").Select()
